# Edit script for "Bazy danych - 3 zajecia.pptx"
#
# Changes applied (per the authoritative OOXML diff):
#   1. Delete slide 17 ("Sprawdz siebie!", sldId 394) - delete the
#      higher-numbered slide first so earlier indices stay stable.
#   2. Delete slide 2 ("WEJSCIOWKA", sldId 395).
#   3. On the slide titled "Podstawy jezyka SQL. Zapytania
#      SELECT-FROM-WHERE. Sortowanie danych." (now at index 2 after the
#      two deletions above), append a blank paragraph and a new
#      paragraph with the text "https://sqliteonline.com/" to the body
#      placeholder that already lists the three topics covered.

$p = $ppt.ActivePresentation

# --- 1 & 2: remove the two slides that were deleted in this revision ---
$p.Slides.Item(17).Delete()
$p.Slides.Item(2).Delete()

# --- 3: append the new link paragraph to the content placeholder ---
$s = $p.Slides.Item(2)
$targetShape = $null
foreach ($sh in $s.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "*Sortowanie danych przy u?yciu ORDER BY*") {
            $targetShape = $sh
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$null = $tr.InsertAfter("`r`rhttps://sqliteonline.com/")
